# Update the data table on the active worksheet:
#  - B6 was stored as text "20"; convert it to a real number.
#  - Append two new rows of submitted feedback data (rows 7 and 8),
#    extending the used range from A1:C6 to A1:C8.
#    Row 8's Age/Comments values ("2" and "1") are kept as text, matching
#    how they were originally submitted, not as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B6: was text "20" -> numeric 20
$ws.Range("B6").Value = 20

# New row 7: gb / 22 / a
$ws.Range("A7").Value = "gb"
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = "a"

# New row 8: j / "2" / "1"  (B8 and C8 must remain text, not numbers)
$ws.Range("A8").Value = "j"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2"
$ws.Range("B8").ClearFormats()

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "1"
$ws.Range("C8").ClearFormats()

Write-Host "Updated sheet: B6 -> number, added rows 7-8 (dimension now A1:C8)"
